$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (team_specific_matrix sheet)
# reflecting additional simulated games and recomputed probabilities.

# Row 2
$ws.Range("B2").Value = 0.1939058171745152
$ws.Range("C2").Value = 0.5290858725761773
$ws.Range("J2").Value = 0.01939058171745152
$ws.Range("P2").Value = 0.1440443213296399
$ws.Range("S2").Value = 0.1135734072022161

# Row 3
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.0101010101010101
$ws.Range("J3").Value = 0.0202020202020202
$ws.Range("P3").Value = 0.8080808080808081
$ws.Range("S3").Value = 0.1464646464646465

# Row 4
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6274509803921569
$ws.Range("S4").Value = 0.3137254901960784

# Row 6
$ws.Range("B6").Value = 0.0660377358490566
$ws.Range("D6").Value = 0.01415094339622642
$ws.Range("E6").Value = 0.004716981132075472
$ws.Range("F6").Value = 0.06132075471698113
$ws.Range("J6").Value = 0.2641509433962264
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1320754716981132
$ws.Range("R6").Value = 0.08962264150943396
$ws.Range("S6").Value = 0.3490566037735849

# Row 7
$ws.Range("B7").Value = 0.1878453038674033
$ws.Range("D7").Value = 0.01657458563535912
$ws.Range("F7").Value = 0.04419889502762431
$ws.Range("J7").Value = 0.1491712707182321
$ws.Range("O7").Value = 0.02762430939226519
$ws.Range("Q7").Value = 0.1657458563535912
$ws.Range("R7").Value = 0.08839779005524862
$ws.Range("S7").Value = 0.3204419889502763

# Row 8
$ws.Range("B8").Value = 0.1422413793103448
$ws.Range("D8").Value = 0.02586206896551724
$ws.Range("F8").Value = 0.07112068965517242
$ws.Range("J8").Value = 0.1142241379310345
$ws.Range("O8").Value = 0.02586206896551724
$ws.Range("Q8").Value = 0.1831896551724138
$ws.Range("R8").Value = 0.0668103448275862
$ws.Range("S8").Value = 0.3706896551724138

# Row 9
$ws.Range("B9").Value = 0.1376146788990826
$ws.Range("D9").Value = 0.02752293577981652
$ws.Range("F9").Value = 0.04587155963302753
$ws.Range("J9").Value = 0.1192660550458716
$ws.Range("O9").Value = 0.009174311926605505
$ws.Range("Q9").Value = 0.1376146788990826
$ws.Range("R9").Value = 0.09174311926605505
$ws.Range("S9").Value = 0.4311926605504587

# Row 10
$ws.Range("B10").Value = 0.1356382978723404
$ws.Range("D10").Value = 0.02659574468085106
$ws.Range("E10").Value = 0.001773049645390071
$ws.Range("F10").Value = 0.07801418439716312
$ws.Range("J10").Value = 0.1205673758865248
$ws.Range("O10").Value = 0.01684397163120567
$ws.Range("Q10").Value = 0.1861702127659574
$ws.Range("R10").Value = 0.06914893617021277
$ws.Range("S10").Value = 0.3652482269503546

# Row 11
$ws.Range("G11").Value = 0.1407407407407407
$ws.Range("J11").Value = 0.0962962962962963
$ws.Range("K11").Value = 0.1888888888888889
$ws.Range("L11").Value = 0.5592592592592592
$ws.Range("S11").Value = 0.01481481481481482

# Row 12
$ws.Range("G12").Value = 0.7278481012658228
$ws.Range("J12").Value = 0.1962025316455696
$ws.Range("K12").Value = 0.0189873417721519
$ws.Range("L12").Value = 0.02531645569620253
$ws.Range("S12").Value = 0.03164556962025317

# Row 13
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2553191489361702
$ws.Range("S13").Value = 0.0851063829787234

# Row 15
$ws.Range("F15").Value = 0.004878048780487805
$ws.Range("H15").Value = 0.2341463414634146
$ws.Range("I15").Value = 0.02926829268292683
$ws.Range("J15").Value = 0.3024390243902439
$ws.Range("K15").Value = 0.04878048780487805
$ws.Range("M15").Value = 0.01463414634146342
$ws.Range("O15").Value = 0.08292682926829269
$ws.Range("S15").Value = 0.2829268292682927

# Row 16
$ws.Range("F16").Value = 0.02192982456140351
$ws.Range("H16").Value = 0.1885964912280702
$ws.Range("I16").Value = 0.05263157894736842
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.09210526315789473
$ws.Range("M16").Value = 0.02192982456140351
$ws.Range("O16").Value = 0.04824561403508772
$ws.Range("S16").Value = 0.1578947368421053

# Row 17
$ws.Range("F17").Value = 0.02472527472527472
$ws.Range("H17").Value = 0.2087912087912088
$ws.Range("I17").Value = 0.0576923076923077
$ws.Range("J17").Value = 0.3846153846153846
$ws.Range("K17").Value = 0.09340659340659341
$ws.Range("M17").Value = 0.01373626373626374
$ws.Range("O17").Value = 0.06593406593406594
$ws.Range("S17").Value = 0.1510989010989011

# Row 18
$ws.Range("F18").Value = 0.006535947712418301
$ws.Range("H18").Value = 0.2026143790849673
$ws.Range("I18").Value = 0.0457516339869281
$ws.Range("J18").Value = 0.457516339869281
$ws.Range("K18").Value = 0.1045751633986928
$ws.Range("M18").Value = 0.0196078431372549
$ws.Range("O18").Value = 0.03267973856209151
$ws.Range("S18").Value = 0.130718954248366

# Row 19
$ws.Range("F19").Value = 0.02293967714528462
$ws.Range("H19").Value = 0.2234494477485132
$ws.Range("I19").Value = 0.05437553101104503
$ws.Range("J19").Value = 0.35089209855565
$ws.Range("K19").Value = 0.1104502973661852
$ws.Range("M19").Value = 0.02888700084961767
$ws.Range("N19").Value = 0.002548853016142736
$ws.Range("O19").Value = 0.06881903143585387
$ws.Range("S19").Value = 0.1376380628717077
